# Add "Employee ID" column (K) to the Sales sheet, cycling employee ids 1011-1016
$wb = $excel.ActiveWorkbook
$sales = $wb.Worksheets.Item("Sales")

$sales.Range("K1").Value = "Employee ID"

$ids = @(1011, 1012, 1013, 1014, 1015, 1016)
for ($r = 2; $r -le 201; $r++) {
    $idx = ($r - 2) % 6
    $sales.Cells.Item($r, 11).Value = $ids[$idx]
}

# Create the new "Employee" worksheet at the end of the workbook
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$emp = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$emp.Name = "Employee"

# Header row
$emp.Range("A1").Value = "Employee ID"
$emp.Range("B1").Value = "Manager ID"
$emp.Range("C1").Value = "Employee"
$emp.Range("D1").Value = "Manager"

$headerRange = $emp.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.Font.Name = "Calibri"
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Data rows
$emp.Range("A2").Value = 1010
$emp.Range("C2").Value = "Roy F"

$emp.Range("A3").Value = 1011
$emp.Range("B3").Value = 1010
$emp.Range("C3").Value = "Pam H"
$emp.Range("D3").Value = "Roy F"

$emp.Range("A4").Value = 1012
$emp.Range("B4").Value = 1010
$emp.Range("C4").Value = "Guy L"
$emp.Range("D4").Value = "Roy F"

$emp.Range("A5").Value = 1013
$emp.Range("B5").Value = 1011
$emp.Range("C5").Value = "Roger M"
$emp.Range("D5").Value = "Pam H"

$emp.Range("A6").Value = 1014
$emp.Range("B6").Value = 1011
$emp.Range("C6").Value = "Kaylie S"
$emp.Range("D6").Value = "Pam H"

$emp.Range("A7").Value = 1015
$emp.Range("B7").Value = 1012
$emp.Range("C7").Value = "Mike O"
$emp.Range("D7").Value = "Guy L"

$emp.Range("A8").Value = 1016
$emp.Range("B8").Value = 1012
$emp.Range("C8").Value = "Rudy Q"
$emp.Range("D8").Value = "Guy L"

# Column widths to match authored layout
$emp.Columns.Item(1).ColumnWidth = 17
$emp.Columns.Item(2).ColumnWidth = 10.77734375
$emp.Columns.Item(3).ColumnWidth = 12.88671875
$emp.Columns.Item(4).ColumnWidth = 11.5546875
